$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$ws.Rows.Item(10).Delete()
$ws.Range("Z50").Select()
